$d = $word.ActiveDocument

$old = "sudo pacman -S inter-font ttf-hack ttf-fira-sans ttf-fira-mono ttf-ibm-plex noto-fonts noto-fonts-emoji ttf-dejavu ttf-liberation ttf-carlito ttf-caladea terminus-font ttf-material-icons ttf-material-symbols-variable ttf-meslo-nerd"

$rng = $d.Content
$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($rng.Find.Found) {
    $start = $rng.Start
    $end = $rng.End

    $part1Text = "sudo pacman -S inter-font noto-fonts noto-fonts-emoji terminus-font ttf-caladea ttf-carlito ttf-dejavu "
    $part2Text = "ttf-fira-code "
    $part3Text = "ttf-fira-mono ttf-fira-sans ttf-hack ttf-ibm-plex ttf-liberation ttf-libertinus ttf-material-icons ttf-material-symbols-variable ttf-meslo-nerd"

    # Replace the whole run's text with the new (reordered + augmented) text.
    $whole = $d.Range($start, $end)
    $whole.Text = $part1Text + $part2Text + $part3Text

    $p1end = $start + $part1Text.Length
    $p2end = $p1end + $part2Text.Length
    $p3end = $p2end + $part3Text.Length

    # Force the middle and trailing segments into their own runs (matching
    # the target's run split) even though their final formatting is the
    # same as the first run: toggle Bold on/off to break run-merging.
    $r2 = $d.Range($p1end, $p2end)
    $r2.Bold = 1
    $r2.Bold = 0

    $r3 = $d.Range($p2end, $p3end)
    $r3.Bold = 1
    $r3.Bold = 0
}

Write-Output $rng.Find.Found
